$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell C10 ("From" value for rule R30, i.e. R30.From) from 18 to 1
$ws.Range("C10").Value = 1
